$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medical Devices")

# Insert 4 new blank rows above row 20 (Dexcom..Sartorius block shifts down by 4)
$ws.Rows.Item(20).Resize(4).Insert()

# New row 20: Idexx Labs
$ws.Range("A20").Value = "x"
$ws.Range("B20").Value = "Idexx Labs"
$ws.Range("C20").Value = "IDXX"

# New row 21: Dexcom (re-entered; old data now sits at row 24 after the insert shift)
$ws.Range("A21").Value = "x"
$ws.Range("B21").Value = "Dexcom"
$ws.Range("C21").Value = "DXCM"

# New row 22: Sartorius (moved up from its old location, now at row 38 after the shift)
$ws.Range("A22").Value = "x"
$ws.Range("B22").Value = "Sartorius"
$ws.Range("C22").Value = "SRT GR"

# New row 23: Fujifilm (brand new company)
$ws.Range("A23").Value = "x"
$ws.Range("B23").Value = "Fujifilm"
$ws.Range("C23").Value = "4901 JP"

# New row 24: Mettler-Toledo (overwrite the shifted-down old Dexcom row)
$ws.Range("A24").Value = "x"
$ws.Range("B24").Value = "Mettler-Toledo"
$ws.Range("C24").Value = "MTD"

# Row 25 already holds the shifted Olympus row; make sure the "x" marker is present
$ws.Range("A25").Value = "x"

# Remove the now-duplicated old Sartorius row (shifted down to row 38)
$ws.Rows.Item(38).Delete()

# Update selection to match the final state
$ws.Range("B25").Select()
